$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("G4").Value = 3.05
$ws.Range("H4").Value = 2.67
$ws.Range("I4").Value = 2.65
$ws.Range("J4").Value = 3.6
$ws.Range("K4").Value = 1.9
$ws.Range("X4").Value = 2.4
$ws.Range("AA4").Value = 7.6
$ws.Range("AB4").Value = 15
$ws.Range("AD4").Value = 40
$ws.Range("AF4").Value = 40
$ws.Range("AI4").Value = 14.5
$ws.Range("AJ4").Value = 80
$ws.Range("AL4").Value = 6.6
$ws.Range("AM4").Value = 12
$ws.Range("AO4").Value = 32

# Row 7
$ws.Range("G7").Value = 1.62
$ws.Range("I7").Value = 5.5
$ws.Range("AA7").Value = 5.5
$ws.Range("AN7").Value = 19
$ws.Range("AO7").Value = 67

# Row 8
$ws.Range("I8").Value = 12
$ws.Range("K8").Value = 2.5
$ws.Range("L8").Value = 12
$ws.Range("W8").Value = 1.36
$ws.Range("X8").Value = 3
$ws.Range("Y8").Value = 2.75
$ws.Range("Z8").Value = 1.4
$ws.Range("AA8").Value = 5.5
$ws.Range("AB8").Value = 5
$ws.Range("AC8").Value = 10
$ws.Range("AE8").Value = 13
$ws.Range("AG8").Value = 10
$ws.Range("AH8").Value = 12
$ws.Range("AI8").Value = 34
$ws.Range("AJ8").Value = 151
$ws.Range("AN8").Value = 34
$ws.Range("AO8").Value = 201
$ws.Range("AP8").Value = 101
$ws.Range("AQ8").Value = 101

# Row 9
$ws.Range("N9").Value = 10
$ws.Range("Q9").Value = 2.05
$ws.Range("W9").Value = 1.4
$ws.Range("X9").Value = 2.75

# Row 10
$ws.Range("G10").Value = 2.3
$ws.Range("I10").Value = 3.1
$ws.Range("U10").Value = 4
$ws.Range("V10").Value = 1.25
$ws.Range("Y10").Value = 1.91
$ws.Range("Z10").Value = 1.8
$ws.Range("AB10").Value = 10
$ws.Range("AD10").Value = 21
$ws.Range("AI10").Value = 17
$ws.Range("AN10").Value = 12
$ws.Range("AP10").Value = 29

# Row 11
$ws.Range("G11").Value = 2.45
$ws.Range("I11").Value = 2.9
$ws.Range("L11").Value = 3.75
$ws.Range("O11").Value = 1.44
$ws.Range("P11").Value = 2.75
$ws.Range("Q11").Value = 2.38
$ws.Range("R11").Value = 1.57
$ws.Range("U11").Value = 4.5
$ws.Range("V11").Value = 1.2
$ws.Range("Y11").Value = 1.91
$ws.Range("Z11").Value = 1.8
$ws.Range("AA11").Value = 7
$ws.Range("AB11").Value = 11
$ws.Range("AG11").Value = 7.5
$ws.Range("AI11").Value = 17
$ws.Range("AK11").Value = 401
$ws.Range("AP11").Value = 26
$ws.Range("AQ11").Value = 41
